# Edit script: expand "Biblioteca de Jogos" heading and flesh out the
# "Biblioteca(s) de Jogos" section with descriptive text plus the new
# "Saga Game Library" subsection (Som / Imagem / Memória / Controle).

$d = $word.ActiveDocument

# ---------------------------------------------------------------------
# 1) "Biblioteca de Jogos" -> "Bibliotecas de Jogos", split across three
#    runs ("Biblioteca" + "s" + " de Jogos") to mirror the authored edit.
# ---------------------------------------------------------------------
$p2 = $d.Paragraphs(2)
$p2StartPt = $d.Range($p2.Range.Start, $p2.Range.Start)
$xml1 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Biblioteca</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>s</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> de Jogos</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p2StartPt.InsertXML($xml1)

$oldHeading = $d.Range($p2.Range.Start, $p2.Range.End)
$oldHeading.Find.Execute("Biblioteca de Jogos")
$oldHeading.Delete()

# ---------------------------------------------------------------------
# 2) The paragraph that used to hold a single space becomes the
#    introductory paragraph about game-engine tooling.
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4Start = $p4.Range.Start
$p4StartPt = $d.Range($p4Start, $p4Start)
$xml2 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">O nível dessas ferramentas varia: algumas </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>engines</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> se limitam a códigos, ou seja, constantes, variáveis, funções e classes relacionadas, mas outras contam com interfaces gráficas que possibilitam o desenvolvimento de um jogo sem programação. De qualquer forma, uma game </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>engine</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve"> precisa proporcionar</w:t></w:r><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">, entre outras funcionalidades, ferramentas para manipular sons, imagens (elementos, texto, imagens, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">), memória (dados) e controle (teclado, mouse, </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>etc</w:t></w:r><w:proofErr w:type="spellEnd"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>).</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4StartPt.InsertXML($xml2)

$insertedLen = 452
$oldSpace = $d.Range($p4Start + $insertedLen, $p4.Range.End)
$oldSpace.Delete()

# ---------------------------------------------------------------------
# 3) Append the new "Saga Game Library" subsection with its Som / Imagem
#    / Memória / Controle sub-headings (and the blank spacer paragraphs
#    that separate them).
# ---------------------------------------------------------------------
$p4 = $d.Paragraphs(4)
$p4EndPt = $d.Range($p4.Range.End, $p4.Range.End)
$xml3 = '<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:pPr><w:pStyle w:val="Ttulo2"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t xml:space="preserve">Saga Game </w:t></w:r><w:proofErr w:type="spellStart"/><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Library</w:t></w:r><w:proofErr w:type="spellEnd"/></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Som</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Imagem</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Memória</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:pStyle w:val="Ttulo3"/><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr><w:r><w:rPr><w:lang w:val="pt-BR"/></w:rPr><w:t>Controle</w:t></w:r></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p><w:p><w:pPr><w:rPr><w:lang w:val="pt-BR"/></w:rPr></w:pPr></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>'
$p4EndPt.InsertXML($xml3)

Write-Output ("Paragraphs.Count=" + $d.Paragraphs.Count)
